$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New last row for the table: index 24 (1-based header at row 1, data starts
# at row 2, so data row 24 -> sheet row 25).
$newRow = 25
$prevRow = 24

# Copy formatting only from the previous data row so the new row picks up
# the same styles (bordered index column A, date/time number format on E)
# that the rest of the table uses.
$ws.Range("A24:V24").Copy()
$ws.Range("A25:V25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item($newRow, 1).Value = 24
$ws.Cells.Item($newRow, 2).Value = "gibraltar"
$ws.Cells.Item($newRow, 3).Value = "national-league"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45233.875
$ws.Cells.Item($newRow, 6).Value = "Magpies"
$ws.Cells.Item($newRow, 7).Value = 2
$ws.Cells.Item($newRow, 8).Value = "St Josephs"
$ws.Cells.Item($newRow, 9).Value = 4

$ws.Cells.Item($newRow, 10).Value = 2.27
$ws.Cells.Item($newRow, 11).Value = "03/11/2023 11:47"
$ws.Cells.Item($newRow, 12).Value = 2.46
$ws.Cells.Item($newRow, 13).Value = "03/11/2023 20:42"
$ws.Cells.Item($newRow, 14).Value = 3.51
$ws.Cells.Item($newRow, 15).Value = "03/11/2023 11:47"
$ws.Cells.Item($newRow, 16).Value = 3.29
$ws.Cells.Item($newRow, 17).Value = "03/11/2023 20:42"
$ws.Cells.Item($newRow, 18).Value = 2.58
$ws.Cells.Item($newRow, 19).Value = "03/11/2023 11:47"
$ws.Cells.Item($newRow, 20).Value = 2.54
$ws.Cells.Item($newRow, 21).Value = "03/11/2023 20:42"

$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/gibraltar/national-league/magpies-st-josephs/vw6GDCFE/"
